$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Harian - Table")

# The sheet had a BMKG logo picture anchored at A1 - remove it as part of
# this cleanup/preprocessing pass.
if ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete()
}

# A1 on the original sheet was always an empty placeholder cell - clear it
# explicitly so it no longer round-trips as a stray empty <c/>.
$ws.Range("A1").ClearContents()

# Copy the whole daily-data table (header row + the 32 data rows) from the
# original sheet into a brand new sheet, pasted starting at A1, keeping the
# new sheet positioned right after the source sheet.
$newSheet = $wb.Worksheets.Add($null, $ws)
$ws.Range("A9:K40").Copy($newSheet.Range("A1"))

# Restore/update the on-screen selections to match the post-edit state:
# the source sheet keeps the just-copied range selected, while the new
# sheet (now the active tab) has its whole used range selected.
$ws.Range("A9:K40").Select()
$newSheet.Range("A1:K32").Select()
$newSheet.Activate()
